# TournRPG-31: add "level up" and "full HP heal" entries to the message
# table, matching the shape of the existing rows (id / msg / color).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("message")

# Row 18 is the last existing data row; clone its formatting (cell
# styles, borders, fill, font) down onto the two new rows 19-20 so they
# look identical to the rest of the table.
$ws1.Range("A18:C18").Copy()
$ws1.Range("A19:C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the row height used by every other data row in the sheet.
$ws1.Rows.Item(19).RowHeight = 20
$ws1.Rows.Item(20).RowHeight = 20

# Column A is the running id, computed the same way as the rows above it.
$ws1.Range("A19").Formula = "=ROW()-2"
$ws1.Range("A20").Formula = "=ROW()-2"

# New message row 17: "<val1>はレベルアップした" (level up), shown in orange.
$ws1.Range("B19").Value = "<val1>はレベルアップした"
$ws1.Range("C19").Value = "orange"

# New message row 18: "<val1>はHPが全回復した" (full HP recovery), shown in blue.
$ws1.Range("B20").Value = "<val1>はHPが全回復した"
$ws1.Range("C20").Value = "blue"
